{"js": "// Replace the 25 two-digit-by-two-digit multiplication prompts in the\n// table with their new values, matching the canonical OOXML diff.\n// Each old prompt text is unique in the document, so a direct\n// search-and-replace on context.document.body is unambiguous.\n\nconst replacements = [\n  [\"21\u00d772=\", \"22\u00d786=\"],\n  [\"27\u00d792=\", \"97\u00d791=\"],\n  [\"77\u00d757=\", \"93\u00d755=\"],\n  [\"40\u00d783=\", \"44\u00d755=\"],\n  [\"53\u00d755=\", \"84\u00d738=\"],\n  [\"84\u00d777=\", \"40\u00d764=\"],\n  [\"14\u00d785=\", \"85\u00d788=\"],\n  [\"37\u00d773=\", \"11\u00d736=\"],\n  [\"89\u00d763=\", \"13\u00d790=\"],\n  [\"41\u00d775=\", \"52\u00d738=\"],\n  [\"27\u00d738=\", \"50\u00d790=\"],\n  [\"11\u00d766=\", \"90\u00d751=\"],\n  [\"59\u00d751=\", \"98\u00d767=\"],\n  [\"53\u00d715=\", \"87\u00d769=\"],\n  [\"83\u00d743=\", \"15\u00d777=\"],\n  [\"97\u00d760=\", \"24\u00d768=\"],\n  [\"31\u00d744=\", \"92\u00d724=\"],\n  [\"91\u00d786=\", \"94\u00d716=\"],\n  [\"55\u00d794=\", \"24\u00d795=\"],\n  [\"31\u00d716=\", \"71\u00d724=\"],\n  [\"53\u00d724=\", \"22\u00d756=\"],\n  [\"91\u00d740=\", \"89\u00d746=\"],\n  [\"47\u00d748=\", \"99\u00d729=\"],\n  [\"53\u00d725=\", \"54\u00d727=\"],\n  [\"23\u00d765=\", \"94\u00d743=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-by-two-digit multiplication prompts in the\n# table with their new values, matching the canonical OOXML diff.\n# Each old prompt text is unique in the document, so Find/Replace on\n# the whole document Content range is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"21\u00d772=\", \"22\u00d786=\"),\n  @(\"27\u00d792=\", \"97\u00d791=\"),\n  @(\"77\u00d757=\", \"93\u00d755=\"),\n  @(\"40\u00d783=\", \"44\u00d755=\"),\n  @(\"53\u00d755=\", \"84\u00d738=\"),\n  @(\"84\u00d777=\", \"40\u00d764=\"),\n  @(\"14\u00d785=\", \"85\u00d788=\"),\n  @(\"37\u00d773=\", \"11\u00d736=\"),\n  @(\"89\u00d763=\", \"13\u00d790=\"),\n  @(\"41\u00d775=\", \"52\u00d738=\"),\n  @(\"27\u00d738=\", \"50\u00d790=\"),\n  @(\"11\u00d766=\", \"90\u00d751=\"),\n  @(\"59\u00d751=\", \"98\u00d767=\"),\n  @(\"53\u00d715=\", \"87\u00d769=\"),\n  @(\"83\u00d743=\", \"15\u00d777=\"),\n  @(\"97\u00d760=\", \"24\u00d768=\"),\n  @(\"31\u00d744=\", \"92\u00d724=\"),\n  @(\"91\u00d786=\", \"94\u00d716=\"),\n  @(\"55\u00d794=\", \"24\u00d795=\"),\n  @(\"31\u00d716=\", \"71\u00d724=\"),\n  @(\"53\u00d724=\", \"22\u00d756=\"),\n  @(\"91\u00d740=\", \"89\u00d746=\"),\n  @(\"47\u00d748=\", \"99\u00d729=\"),\n  @(\"53\u00d725=\", \"54\u00d727=\"),\n  @(\"23\u00d765=\", \"94\u00d743=\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Text = $old\n  $rng.Find.Replacement.Text = $new\n  $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n}\n"}
